# Active_Outages.xlsx refresh - 6/19/2025, 9:13:27 AM
# 1) Elapsed Duration(Hrs) values advance as the report is regenerated.
# 2) The JED0125 outage (R1 sheet, row 6) has been resolved (0:00:00 duration)
#    and is dropped from the active-outages list.

$wb = $excel.ActiveWorkbook

# --- R1 sheet ---
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3946:27:33"
$ws.Range("G3").Value = "86:00:11"
$ws.Range("G4").Value = "109:00:11"
$ws.Rows.Item(6).Delete()

# --- R2 sheet ---
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12127:51:14"
$ws.Range("G3").Value = "3257:34:43"
$ws.Range("G4").Value = "495:46:17"

# --- R4 sheet ---
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2973:41:03"
$ws.Range("G3").Value = "200:53:18"
$ws.Range("G4").Value = "89:05:43"
$ws.Range("G5").Value = "86:43:16"

# --- R5 sheet ---
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "447:40:02"

# --- R6 sheet ---
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "88:12:20"
